$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-valued columns (D: Price, E: Volume) to remain text even when
# numeric-looking, matching the source inlineStr cells; reset the style
# afterwards so no stray NumberFormat/quotePrefix style lingers on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '47.895.42'
Set-TextValue $ws.Range("E2") '  -0.99%  '
Set-TextValue $ws.Range("D3") '2.485.63'
Set-TextValue $ws.Range("E3") '  -1.48%  '
Set-TextValue $ws.Range("D4") '0.998'
Set-TextValue $ws.Range("E4") '  -0.22%  '
Set-TextValue $ws.Range("D5") '317.27'
Set-TextValue $ws.Range("E5") '  -1.59%  '
Set-TextValue $ws.Range("D6") '104.28'
Set-TextValue $ws.Range("E6") '  -5.01%  '
Set-TextValue $ws.Range("D7") '0.518'
Set-TextValue $ws.Range("E7") '  -2.61%  '
Set-TextValue $ws.Range("D8") '0.999'
Set-TextValue $ws.Range("E8") '  -0.08%  '
Set-TextValue $ws.Range("D9") '0.535'
Set-TextValue $ws.Range("E9") '  -2.98%  '
Set-TextValue $ws.Range("D10") '38.77'
Set-TextValue $ws.Range("E10") '  -4.42%  '
Set-TextValue $ws.Range("D11") '20.32'
Set-TextValue $ws.Range("E11") '  -0.48%  '
Set-TextValue $ws.Range("E12") '  -2.97%  '
Set-TextValue $ws.Range("E13") '  +0.54%  '
Set-TextValue $ws.Range("D14") '7.02'
Set-TextValue $ws.Range("E14") '  -3.52%  '
Set-TextValue $ws.Range("D15") '2.870.77'
Set-TextValue $ws.Range("E15") '  -1.72%  '
Set-TextValue $ws.Range("D16") '2.527.00'
Set-TextValue $ws.Range("E16") '  +0.38%  '
Set-TextValue $ws.Range("D17") '0.824'
Set-TextValue $ws.Range("E17") '  -3.50%  '
Set-TextValue $ws.Range("D18") '47.769.47'
Set-TextValue $ws.Range("E18") '  -0.88%  '
Set-TextValue $ws.Range("D19") '12.67'
Set-TextValue $ws.Range("E19") '  -4.98%  '
Set-TextValue $ws.Range("D20") '2.88'
Set-TextValue $ws.Range("E20") '  +7.14%  '
Set-TextValue $ws.Range("D21") '6.51'
Set-TextValue $ws.Range("E21") '  -1.75%  '
Set-TextValue $ws.Range("E22") '  -2.55%  '
Set-TextValue $ws.Range("D23") '277.94'
Set-TextValue $ws.Range("E23") '  +1.65%  '
Set-TextValue $ws.Range("D24") '70.57'
Set-TextValue $ws.Range("E24") '  -2.06%  '
Set-TextValue $ws.Range("D25") '2.48'
Set-TextValue $ws.Range("E25") '  -4.16%  '
Set-TextValue $ws.Range("E26") '  -0.15%  '
Set-TextValue $ws.Range("D27") '25.60'
Set-TextValue $ws.Range("E27") '  -1.58%  '
Set-TextValue $ws.Range("E28") '  -7.25%  '
Set-TextValue $ws.Range("D29") '9.56'
Set-TextValue $ws.Range("E29") '  -5.64%  '
Set-TextValue $ws.Range("D30") '0.138'
Set-TextValue $ws.Range("E30") '  -5.07%  '
Set-TextValue $ws.Range("E31") '  -3.82%  '
Set-TextValue $ws.Range("D32") '48.96'
Set-TextValue $ws.Range("E32") '  -1.40%  '
Set-TextValue $ws.Range("E33") '  -0.32%  '
Set-TextValue $ws.Range("D35") '5.24'
Set-TextValue $ws.Range("E35") '  -2.83%  '
Set-TextValue $ws.Range("D36") '0.0767'
Set-TextValue $ws.Range("E36") '  -2.52%  '
Set-TextValue $ws.Range("E37") '  -2.77%  '
Set-TextValue $ws.Range("D38") '4.49'
Set-TextValue $ws.Range("E38") '  -4.29%  '
Set-TextValue $ws.Range("E39") '  -5.18%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D40") '0.110'
Set-TextValue $ws.Range("E40") '  -1.48%  '
$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '2.21'
Set-TextValue $ws.Range("E41") '  -0.36%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D42") '119.92'
Set-TextValue $ws.Range("E42") '  -2.44%  '
Set-TextValue $ws.Range("D43") '21.36'
Set-TextValue $ws.Range("E43") '  -3.92%  '
Set-TextValue $ws.Range("D44") '0.0298'
Set-TextValue $ws.Range("E44") '  -1.35%  '
Set-TextValue $ws.Range("D45") '1.985.11'
Set-TextValue $ws.Range("E45") '  -2.16%  '
Set-TextValue $ws.Range("E46") '  -1.45%  '
Set-TextValue $ws.Range("D47") '1.91'
Set-TextValue $ws.Range("E47") '  +0.87%  '
Set-TextValue $ws.Range("E48") '  +0.04%  '
Set-TextValue $ws.Range("D49") '8.89'
Set-TextValue $ws.Range("E49") '  -2.57%  '
Set-TextValue $ws.Range("E50") '  -2.95%  '
Set-TextValue $ws.Range("D51") '78.89'
Set-TextValue $ws.Range("E51") '  -1.18%  '
